$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 442.8889
$ws.Range("J4").Value = 695.4
$ws.Range("L4").Value = 695.4
$ws.Range("N4").Value = -923.4
$ws.Range("H5").Value = 55
$ws.Range("I5").Value = 55
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 55
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 60
$ws.Range("N5").Value = ""
$ws.Range("H12").Value = 407.42856
$ws.Range("I12").Value = 205
$ws.Range("J12").Value = 677.3333
$ws.Range("K12").Value = 205
$ws.Range("L12").Value = 677.3333
$ws.Range("M12").Value = -35
$ws.Range("N12").Value = -1017.3333
$ws.Range("H21").Value = 13345.333
$ws.Range("I21").Value = 13345.333
$ws.Range("K21").Value = 13345.333
$ws.Range("M21").Value = -12877.333
$ws.Range("H23").Value = 13345.333
$ws.Range("I23").Value = 13345.333
$ws.Range("K23").Value = 13345.333
$ws.Range("M23").Value = -13111.333
$ws.Range("H51").Value = 6814.5
$ws.Range("I51").Value = 7471.75
$ws.Range("J51").Value = 5500
$ws.Range("K51").Value = 7471.75
$ws.Range("L51").Value = 5500
$ws.Range("M51").Value = -6987.75
$ws.Range("N51").Value = -6468
$ws.Range("H62").Value = 7493
$ws.Range("I62").Value = 4543.091
$ws.Range("K62").Value = 4543.091
$ws.Range("M62").Value = -3919.091
$ws.Range("H64").Value = 6197.6924
$ws.Range("I64").Value = 4224.2856
$ws.Range("J64").Value = 8500
$ws.Range("K64").Value = 4224.2856
$ws.Range("L64").Value = 8500
$ws.Range("M64").Value = -3976.2856
$ws.Range("N64").Value = -8996
$ws.Range("H65").Value = 7493
$ws.Range("I65").Value = 4543.091
$ws.Range("K65").Value = 22715.455
$ws.Range("M65").Value = -19595.455
$ws.Range("H67").Value = 6197.6924
$ws.Range("I67").Value = 4224.2856
$ws.Range("J67").Value = 8500
$ws.Range("K67").Value = 4224.2856
$ws.Range("L67").Value = 8500
$ws.Range("M67").Value = -3366.2856
$ws.Range("N67").Value = -10216
$ws.Range("H135").Value = 2397.5
$ws.Range("I135").Value = 1684.3334
$ws.Range("J135").Value = 3110.6667
$ws.Range("K135").Value = 15159.0006
$ws.Range("L135").Value = 27996.0003
$ws.Range("M135").Value = -12624.0006
$ws.Range("N135").Value = -33066.0003
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5495.0435
$ws.Range("I2").Value = 1467.4375
$ws.Range("J2").Value = 14701
$ws.Range("K2").Value = 1467.4375
$ws.Range("L2").Value = 14701
$ws.Range("M2").Value = -1354.4375
$ws.Range("N2").Value = -14927
$ws.Range("H61").Value = 3008.7273
$ws.Range("I61").Value = 2809.6
$ws.Range("K61").Value = 2809.6
$ws.Range("M61").Value = -2597.6
$ws.Range("H63").Value = 1128.4445
$ws.Range("I63").Value = 879.7143
$ws.Range("J63").Value = 1999
$ws.Range("K63").Value = 879.7143
$ws.Range("L63").Value = 1999
$ws.Range("M63").Value = -193.7143
$ws.Range("N63").Value = -3371
$ws.Range("H66").Value = 1128.4445
$ws.Range("I66").Value = 879.7143
$ws.Range("J66").Value = 1999
$ws.Range("K66").Value = 4398.5715
$ws.Range("L66").Value = 9995
$ws.Range("M66").Value = -966.5715
$ws.Range("N66").Value = -16859
$ws.Range("H102").Value = 3174.5334
$ws.Range("I102").Value = 692.5454999999999
$ws.Range("J102").Value = 10000
$ws.Range("K102").Value = 692.5454999999999
$ws.Range("L102").Value = 10000
$ws.Range("M102").Value = 929.4545000000001
$ws.Range("N102").Value = -13244
$ws.Range("H105").Value = 273342.5
$ws.Range("J105").Value = 273342.5
$ws.Range("L105").Value = 273342.5
$ws.Range("N105").Value = -280330.5
$ws.Range("H116").Value = 5495.0435
$ws.Range("I116").Value = 1467.4375
$ws.Range("J116").Value = 14701
$ws.Range("K116").Value = 1467.4375
$ws.Range("L116").Value = 14701
$ws.Range("M116").Value = 826.5625
$ws.Range("N116").Value = -19289
$ws.Range("H136").Value = 3008.7273
$ws.Range("I136").Value = 2809.6
$ws.Range("K136").Value = 8428.799999999999
$ws.Range("M136").Value = -5878.799999999999
$ws.Range("H138").Value = 100000.5
$ws.Range("J138").Value = 100000.5
$ws.Range("L138").Value = 100000.5
$ws.Range("N138").Value = -110280.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5495.0435
$ws.Range("I3").Value = 1467.4375
$ws.Range("J3").Value = 14701
$ws.Range("K3").Value = 1467.4375
$ws.Range("L3").Value = 14701
$ws.Range("M3").Value = -1353.4375
$ws.Range("N3").Value = -14929
$ws.Range("H22").Value = 195.83333
$ws.Range("I22").Value = 143.75
$ws.Range("J22").Value = 300
$ws.Range("K22").Value = 143.75
$ws.Range("L22").Value = 300
$ws.Range("M22").Value = 29.25
$ws.Range("N22").Value = -646
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").Value = ""
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").Value = ""
$ws.Range("H116").Value = 54392.5
$ws.Range("J116").Value = 54392.5
$ws.Range("L116").Value = 54392.5
$ws.Range("N116").Value = -63570.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 583.375
$ws.Range("I5").Value = 238.14285
$ws.Range("J5").Value = 3000
$ws.Range("K5").Value = 238.14285
$ws.Range("L5").Value = 3000
$ws.Range("M5").Value = -126.14285
$ws.Range("N5").Value = -3224
$ws.Range("H55").Value = 42750
$ws.Range("J55").Value = 42750
$ws.Range("L55").Value = 42750
$ws.Range("N55").Value = -43380
$ws.Range("H58").Value = 3935.4167
$ws.Range("I58").Value = 1321.5714
$ws.Range("K58").Value = 1321.5714
$ws.Range("M58").Value = -1118.5714
$ws.Range("H64").Value = 57750
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 57750
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 57750
$ws.Range("M64").Value = ""
$ws.Range("N64").Value = -58246
$ws.Range("H67").Value = 57750
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 57750
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 57750
$ws.Range("M67").Value = ""
$ws.Range("N67").Value = -59466
$ws.Range("H103").Value = 0
$ws.Range("I103").Value = 0
$ws.Range("K103").Value = 0
$ws.Range("M103").Value = ""
$ws.Range("H132").Value = 2941
$ws.Range("I132").Value = 2228.2222
$ws.Range("J132").Value = 6148.5
$ws.Range("K132").Value = 6684.6666
$ws.Range("L132").Value = 18445.5
$ws.Range("M132").Value = -4154.6666
$ws.Range("N132").Value = -23505.5
$ws.Range("H136").Value = 3935.4167
$ws.Range("I136").Value = 1321.5714
$ws.Range("K136").Value = 3964.7142
$ws.Range("M136").Value = -1414.7142
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 4077.5557
$ws.Range("J55").Value = 4799.7144
$ws.Range("L55").Value = 14399.1432
$ws.Range("N55").Value = -14753.1432
$ws.Range("H68").Value = 700
$ws.Range("J68").Value = 1200
$ws.Range("L68").Value = 3600
$ws.Range("N68").Value = -5222
$ws.Range("H71").Value = 700
$ws.Range("J71").Value = 1200
$ws.Range("L71").Value = 10800
$ws.Range("N71").Value = -18912
$ws.Range("H98").Value = 223.71428
$ws.Range("I98").Value = 147.5
$ws.Range("J98").Value = 254.2
$ws.Range("K98").Value = 442.5
$ws.Range("L98").Value = 762.5999999999999
$ws.Range("M98").Value = 1055.5
$ws.Range("N98").Value = -3758.6
$ws.Range("H129").Value = 2362.4285
$ws.Range("J129").Value = 5383.4
$ws.Range("L129").Value = 16150.2
$ws.Range("N129").Value = -26150.2
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 5950
$ws.Range("J4").Value = 5950
$ws.Range("L4").Value = 5950
$ws.Range("N4").Value = -6174
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").Value = ""
$ws.Range("H9").Value = 249.875
$ws.Range("I9").Value = 230.5
$ws.Range("J9").Value = 308
$ws.Range("K9").Value = 230.5
$ws.Range("L9").Value = 308
$ws.Range("M9").Value = -60.5
$ws.Range("N9").Value = -648
$ws.Range("H107").Value = 361.27274
$ws.Range("I107").Value = 330.44446
$ws.Range("K107").Value = 330.44446
$ws.Range("M107").Value = 1589.55554
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2112.3333
$ws.Range("I16").Value = 535
$ws.Range("K16").Value = 535
$ws.Range("M16").Value = -365
$ws.Range("H93").Value = 1107.7142
$ws.Range("I93").Value = 959
$ws.Range("K93").Value = 959
$ws.Range("M93").Value = 289
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1954.4286
$ws.Range("I100").Value = 1954.4286
$ws.Range("K100").Value = 3908.8572
$ws.Range("M100").Value = -3367.8572
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").Value = ""
